$wb = $excel.ActiveWorkbook

function Set-F {
    param($ws, $row, $val)
    $ws.Cells.Item($row, 6).Value = $val
}

function Set-G {
    param($ws, $row, $val)
    $ws.Cells.Item($row, 7).Value = $val
}

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
Set-F $ws1 3 1200
Set-G $ws1 4 "不可售"
Set-G $ws1 5 75
Set-F $ws1 7 1734
Set-F $ws1 8 427
Set-F $ws1 9 80
Set-F $ws1 10 56
Set-F $ws1 12 275
Set-F $ws1 13 1649
Set-F $ws1 14 311
Set-F $ws1 16 769
Set-F $ws1 17 308
Set-F $ws1 18 646
Set-F $ws1 19 12528
Set-F $ws1 20 12573
Set-F $ws1 22 725
Set-F $ws1 24 286
Set-F $ws1 26 446
Set-F $ws1 27 1952
Set-F $ws1 30 222
Set-F $ws1 31 646

# Sheet 2: 演出 (Show)
$ws2 = $wb.Worksheets.Item("演出")
Set-F $ws2 5 6
Set-F $ws2 6 6

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
Set-F $ws3 2 81
Set-F $ws3 3 134

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
Set-F $ws4 3 81
Set-F $ws4 4 1200
Set-G $ws4 5 "不可售"
Set-G $ws4 6 75
Set-F $ws4 7 134
Set-F $ws4 9 1734
Set-F $ws4 10 427
Set-F $ws4 12 80
Set-F $ws4 13 56
Set-F $ws4 17 275
Set-F $ws4 18 1649
Set-F $ws4 19 311
Set-F $ws4 21 769
Set-F $ws4 22 308
Set-F $ws4 23 6
Set-F $ws4 24 646
Set-F $ws4 25 12528
Set-F $ws4 26 12573
Set-F $ws4 28 725
Set-F $ws4 30 286
Set-F $ws4 32 446
Set-F $ws4 33 6
Set-F $ws4 35 1952
Set-F $ws4 40 222
Set-F $ws4 41 646
